{"js": "// Bold four heading-style lines (name, two job titles, education title) and\n// rename the \"Animation \u2013 Praktikantin\" job-title line to\n// \"Praktikant im Bereich Animation\" (which also becomes bold), matching the\n// source diff. Each run is located via a paragraph-scoped search so that\n// only the exact heading run is touched (not any substring occurrence\n// elsewhere in the document, and not the paragraph mark's own formatting).\n\nconst oldTitle = \"Animation \u2013 Praktikantin\";\nconst newTitle = \"Praktikant im Bereich Animation\";\n\nasync function findParagraphByText(targetText) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.trim() === targetText) {\n      return paragraphs.items[i];\n    }\n  }\n  return null;\n}\n\nasync function boldExactRun(paragraph, text) {\n  const results = paragraph.search(text, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].font.bold = true;\n  }\n  await context.sync();\n}\n\n// 1) \"Holly Dickson\" (name heading) -> bold\nlet para = await findParagraphByText(\"Holly Dickson\");\nawait boldExactRun(para, \"Holly Dickson\");\n\n// 2) \"Senior Animation Designer\" (job title heading) -> bold\npara = await findParagraphByText(\"Senior Animation Designer\");\nawait boldExactRun(para, \"Senior Animation Designer\");\n\n// 3) \"Junior Animation Designer\" (job title heading) -> bold\npara = await findParagraphByText(\"Junior Animation Designer\");\nawait boldExactRun(para, \"Junior Animation Designer\");\n\n// 4) \"Animation \u2013 Praktikantin\" -> replace text, then bold the new text\npara = await findParagraphByText(oldTitle);\nconst oldResults = para.search(oldTitle, { matchCase: true, matchWholeWord: false });\noldResults.load(\"text\");\nawait context.sync();\noldResults.items[0].insertText(newTitle, Word.InsertLocation.replace);\nawait context.sync();\n\npara = await findParagraphByText(newTitle);\nawait boldExactRun(para, newTitle);\n\n// 5) \"Bachelor of Fine Arts in Animation\" (education title heading) -> bold\npara = await findParagraphByText(\"Bachelor of Fine Arts in Animation\");\nawait boldExactRun(para, \"Bachelor of Fine Arts in Animation\");\n", "ps1": "# Bold four heading-style lines (name, two job titles, education title) and\n# rename the \"Animation \u2013 Praktikantin\" job-title line to\n# \"Praktikant im Bereich Animation\" (which also becomes bold), matching the\n# source diff. Each run is located via a paragraph-scoped Find so that only\n# the exact heading run is touched (not any substring occurrence elsewhere\n# in the document, and not the paragraph mark's own formatting).\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphByText($text) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.Trim() -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Set-RunBold($paragraph, $text) {\n    $r = $paragraph.Range.Duplicate\n    $r.Find.Execute($text)\n    $r.Font.Bold = 1\n}\n\n$oldTitle = \"Animation \u2013 Praktikantin\"\n$newTitle = \"Praktikant im Bereich Animation\"\n\n# 1) \"Holly Dickson\" (name heading) -> bold\n$p = Get-ParagraphByText \"Holly Dickson\"\nSet-RunBold $p \"Holly Dickson\"\n\n# 2) \"Senior Animation Designer\" (job title heading) -> bold\n$p = Get-ParagraphByText \"Senior Animation Designer\"\nSet-RunBold $p \"Senior Animation Designer\"\n\n# 3) \"Junior Animation Designer\" (job title heading) -> bold\n$p = Get-ParagraphByText \"Junior Animation Designer\"\nSet-RunBold $p \"Junior Animation Designer\"\n\n# 4) \"Animation \u2013 Praktikantin\" -> replace text, then bold the new text\n$p = Get-ParagraphByText $oldTitle\n$r = $p.Range.Duplicate\n$r.Find.Execute($oldTitle)\n$r.Text = $newTitle\n\n$p = Get-ParagraphByText $newTitle\nSet-RunBold $p $newTitle\n\n# 5) \"Bachelor of Fine Arts in Animation\" (education title heading) -> bold\n$p = Get-ParagraphByText \"Bachelor of Fine Arts in Animation\"\nSet-RunBold $p \"Bachelor of Fine Arts in Animation\"\n"}
